$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("pythonCode")
$ws3 = $wb.Worksheets.Item("Practice Questions")

# --- pythonCode sheet: insert a new TestCaseId column (A) mirroring the
# "Practice Questions" sheet layout, shifting the existing data right ---
$ws2.Columns("A").Insert() | Out-Null

$ws2.Range("A1").Value = "TestCaseId"
$ws2.Range("A2").Value = "TC001"
$ws2.Range("A3").Value = "TC002"
$ws2.Range("A4").Value = "TC003"

$ws2.Columns("A").ColumnWidth = 9.7

# --- Selections / active sheet / active cell bookkeeping ---
$ws2.Range("A1").Select() | Out-Null
$ws3.Range("A2:A4").Select() | Out-Null
$ws2.Activate() | Out-Null
